$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts N..P to O..Q)
$ws.Columns("N:N").Insert()

# The newly inserted column N gets its own width (matches target col width 11)
$ws.Columns("N:N").ColumnWidth = 10.2

# Update the selection / active sheet to match the target workbook view
$ws.Activate()
$ws.Range("R7").Select()
